$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7, shifting existing rows 7-91 down to 8-92.
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with the new data record.
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C7").Value = "Coquimbo"
$ws.Range("D7").Value = 45092
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 100112026
$ws.Range("G7").Value = "Haba"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 1200
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = 13500
$ws.Range("N7").Value = "$/saco 25 kilos"
$ws.Range("O7").Value = "Provincia de Limarí"
$ws.Range("P7").Value = 540
$ws.Range("Q7").Value = 25
$ws.Range("R7").Value = "Hortaliza"
